$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "VALOR MORA" total and the trabajador/periodo counters
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 234008
$ws.Range("C13").Value2 = 4
$ws.Range("F13").Value2 = 5

# ---------------------------------------------------------------------------
# 2) Make room for two new worker rows right after the current last data
#    row (19), pushing the blank rows + signature block down by 2.
# ---------------------------------------------------------------------------
$ws.Range("B20:J21").EntireRow.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Row 21 becomes the new last data row, so it should carry the heavier
# bottom-border look that row 19 used to have ...
$ws.Range("B19:J19").Copy()
$ws.Range("B21:J21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ... while row 19 (no longer last) and the brand new row 20 switch to the
# regular "middle" row look shared by rows 16-18.
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B20:J20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the two new worker rows.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "73201468"
$ws.Range("D20").Value2 = "OWAR DIAZ DE AVILA"
$ws.Range("E20").Value2 = "2509"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1002248999"
$ws.Range("D21").Value2 = "WHAJIR ANTONIO PEREZ ARAUJO"
$ws.Range("E21").Value2 = "2509"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500

# ---------------------------------------------------------------------------
# 4) Cosmetic touch-ups that came along with the edit in the source file:
#    slightly wider columns to fit the new, longer values, and the logo
#    nudged to the right.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.54296875
$ws.Columns.Item(3).ColumnWidth = 16.7265625
$ws.Columns.Item(4).ColumnWidth = 38.36328125
$ws.Columns.Item(5).ColumnWidth = 13.54296875
$ws.Columns.Item(6).ColumnWidth = 10.1796875
$ws.Columns.Item(7).ColumnWidth = 14.36328125
$ws.Columns.Item(8).ColumnWidth = 19.36328125
$ws.Columns.Item(9).ColumnWidth = 18.08984375
$ws.Columns.Item(10).ColumnWidth = 15

$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left + 13.5
Write-Output "Script completed successfully"
